# Updates the cryptocurrency price (D) and 1h volume-change (E) columns
# with refreshed values from the Aug 2 2023 GitHub Actions data pull.
# Price cells that would otherwise be auto-parsed as numbers by Excel
# (losing formatting such as trailing zeros, e.g. "1.000" -> 1) are
# written with a leading apostrophe so they stay plain text, exactly as
# they were stored in the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.621.39"
$ws.Range("E2").Value = "  +2.45%  "

$ws.Range("D3").Value = "1.860.38"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "'244.88"
$ws.Range("E5").Value = "  +1.94%  "

$ws.Range("D6").Value = "'0.6967"
$ws.Range("E6").Value = "  +1.42%  "

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.07706"

$ws.Range("D9").Value = "'0.3061"
$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("D10").Value = "'23.70"
$ws.Range("E10").Value = "  +0.86%  "

$ws.Range("D11").Value = "'0.07748"
$ws.Range("E11").Value = "  -0.45%  "

$ws.Range("D12").Value = "'5.158"
$ws.Range("E12").Value = "  +1.86%  "

$ws.Range("D13").Value = "1.852.11"
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").Value = "'92.20"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").Value = "'0.6925"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").Value = "'6.569"
$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("D17").Value = "29.608.34"
$ws.Range("E17").Value = "  +2.44%  "

$ws.Range("D18").Value = "'0.000008311"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").Value = "2.100.08"
$ws.Range("E19").Value = "  +1.29%  "

$ws.Range("D20").Value = "'241.04"
$ws.Range("E20").Value = "  -0.60%  "

$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("D22").Value = "'0.9996"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").Value = "'7.607"

$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'0.1500"
$ws.Range("E25").Value = "  +2.01%  "

$ws.Range("D26").Value = "'8.922"
$ws.Range("E26").Value = "  +1.86%  "

$ws.Range("D27").Value = "'159.40"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("E28").Value = "  +0.66%  "

$ws.Range("D29").Value = "'1.534"

$ws.Range("E30").Value = "  +1.20%  "

$ws.Range("D31").Value = "'4.188"
$ws.Range("E31").Value = "  +1.74%  "

$ws.Range("D32").Value = "'1.199"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("D34").Value = "'0.7759"
$ws.Range("E34").Value = "  +3.97%  "

$ws.Range("D35").Value = "'1.898"
$ws.Range("E35").Value = "  +4.36%  "

$ws.Range("D36").Value = "'1.153"
$ws.Range("E36").Value = "  +1.07%  "

$ws.Range("D37").Value = "'2.684"
$ws.Range("E37").Value = "  +0.33%  "

$ws.Range("D38").Value = "1.326.70"
$ws.Range("E38").Value = "  +9.59%  "

$ws.Range("E39").Value = "  +1.94%  "

$ws.Range("E40").Value = "  +1.80%  "

$ws.Range("D41").Value = "'0.9738"
$ws.Range("E41").Value = "  +6.37%  "

$ws.Range("D42").Value = "'106.73"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("D43").Value = "'5.825"
$ws.Range("E43").Value = "  +6.75%  "

$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("E45").Value = "  +3.98%  "

$ws.Range("D46").Value = "'9.784"
$ws.Range("E46").Value = "  +2.98%  "

$ws.Range("D47").Value = "2.001.52"
$ws.Range("E47").Value = "  +1.37%  "

$ws.Range("D48").Value = "'0.5215"
$ws.Range("E48").Value = "  +0.95%  "

$ws.Range("D49").Value = "'1.778"
$ws.Range("E49").Value = "  +2.80%  "

$ws.Range("D50").Value = "'63.52"
$ws.Range("E50").Value = "  +0.90%  "

$ws.Range("D51").Value = "'6.968"
$ws.Range("E51").Value = "  +1.24%  "
